# Update cryptocurrency price/volume figures (cryptos list refresh).
# Each cell is forced to Text (leading apostrophe, like typing `27.846.19 in Excel)
# so digit-for-digit formatting (e.g. trailing zeros, grouping dots) is preserved
# exactly as scraped, instead of being auto-parsed into a Number/General value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "27.846.19"
$ws.Range("D3").Value = "'" + "1.906.84"
$ws.Range("E3").Value = "'" + "  -0.11%  "
$ws.Range("E4").Value = "'" + "  -0.29%  "
$ws.Range("D5").Value = "'" + "313.72"
$ws.Range("E5").Value = "'" + "  -0.84%  "
$ws.Range("D6").Value = "'" + "1.004"
$ws.Range("E6").Value = "'" + "  -0.20%  "
$ws.Range("D7").Value = "'" + "0.5025"
$ws.Range("E7").Value = "'" + "  +4.02%  "
$ws.Range("E8").Value = "'" + "  +0.06%  "
$ws.Range("E9").Value = "'" + "  -1.06%  "
$ws.Range("D10").Value = "'" + "0.9077"
$ws.Range("E10").Value = "'" + "  -2.93%  "
$ws.Range("D11").Value = "'" + "20.83"
$ws.Range("E11").Value = "'" + "  +0.08%  "
$ws.Range("D12").Value = "'" + "0.07674"
$ws.Range("E12").Value = "'" + "  -1.29%  "
$ws.Range("D13").Value = "'" + "1.870.32"
$ws.Range("E13").Value = "'" + "  -2.07%  "
$ws.Range("D14").Value = "'" + "5.480"
$ws.Range("E14").Value = "'" + "  -0.42%  "
$ws.Range("D15").Value = "'" + "91.82"
$ws.Range("E15").Value = "'" + "  -0.05%  "
$ws.Range("E16").Value = "'" + "  -0.27%  "
$ws.Range("D17").Value = "'" + "0.000008709"
$ws.Range("E17").Value = "'" + "  -1.41%  "
$ws.Range("E18").Value = "'" + "  -0.23%  "
$ws.Range("D19").Value = "'" + "27.867.37"
$ws.Range("E19").Value = "'" + "  -0.73%  "
$ws.Range("D20").Value = "'" + "14.57"
$ws.Range("D21").Value = "'" + "5.170"
$ws.Range("E21").Value = "'" + "  -0.23%  "
$ws.Range("E22").Value = "'" + "  -0.96%  "
$ws.Range("D23").Value = "'" + "6.578"
$ws.Range("E23").Value = "'" + "  -0.75%  "
$ws.Range("D24").Value = "'" + "154.20"
$ws.Range("E24").Value = "'" + "  -1.00%  "
$ws.Range("D25").Value = "'" + "1.882"
$ws.Range("E25").Value = "'" + "  -2.13%  "
$ws.Range("D26").Value = "'" + "2.215"
$ws.Range("E26").Value = "'" + "  +4.21%  "
$ws.Range("D27").Value = "'" + "18.39"
$ws.Range("E27").Value = "'" + "  -0.66%  "
$ws.Range("D28").Value = "'" + "115.40"
$ws.Range("E28").Value = "'" + "  -1.10%  "
$ws.Range("D29").Value = "'" + "4.902"
$ws.Range("E29").Value = "'" + "  -1.44%  "
$ws.Range("D30").Value = "'" + "0.09012"
$ws.Range("E30").Value = "'" + "  +0.54%  "
$ws.Range("E31").Value = "'" + "  -2.71%  "
$ws.Range("D32").Value = "'" + "1.223"
$ws.Range("E32").Value = "'" + "  -2.25%  "
$ws.Range("D33").Value = "'" + "4.655"
$ws.Range("E33").Value = "'" + "  -0.52%  "
$ws.Range("D34").Value = "'" + "0.7615"
$ws.Range("E34").Value = "'" + "  -1.79%  "
$ws.Range("D35").Value = "'" + "0.02063"
$ws.Range("E35").Value = "'" + "  +0.20%  "
$ws.Range("D36").Value = "'" + "2.492"
$ws.Range("E36").Value = "'" + "  -5.99%  "
$ws.Range("E37").Value = "'" + "  -1.80%  "
$ws.Range("D38").Value = "'" + "0.5528"
$ws.Range("E38").Value = "'" + "  +0.81%  "
$ws.Range("D39").Value = "'" + "3.015"
$ws.Range("E39").Value = "'" + "  +0.67%  "
$ws.Range("D40").Value = "'" + "0.05245"
$ws.Range("E40").Value = "'" + "  -1.14%  "
$ws.Range("D41").Value = "'" + "6.875"
$ws.Range("E41").Value = "'" + "  -2.22%  "
$ws.Range("D42").Value = "'" + "8.465"
$ws.Range("E42").Value = "'" + "  -0.62%  "
$ws.Range("E43").Value = "'" + "  -1.17%  "
$ws.Range("D44").Value = "'" + "110.82"
$ws.Range("E44").Value = "'" + "  +2.42%  "
$ws.Range("D45").Value = "'" + "10.60"
$ws.Range("E45").Value = "'" + "  -1.10%  "
$ws.Range("D46").Value = "'" + "0.4814"
$ws.Range("E46").Value = "'" + "  -0.33%  "
$ws.Range("D47").Value = "'" + "1.004"
$ws.Range("E47").Value = "'" + "  -0.20%  "
$ws.Range("D48").Value = "'" + "1.624"
$ws.Range("E48").Value = "'" + "  -1.61%  "
$ws.Range("D49").Value = "'" + "67.28"
$ws.Range("E49").Value = "'" + "  -0.97%  "
$ws.Range("E50").Value = "'" + "  -0.20%  "
$ws.Range("D51").Value = "'" + "0.9026"
$ws.Range("E51").Value = "'" + "  +0.25%  "
